$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 224, shifting existing rows 224:268 down to 225:269.
$ws.Rows.Item(224).Insert()

$ws.Cells.Item(224, 1).Value = 10
$ws.Cells.Item(224, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(224, 3).Value = "La Araucanía"
$ws.Cells.Item(224, 4).Value = 44995
$ws.Cells.Item(224, 5).Value = 9
$ws.Cells.Item(224, 6).Value = 100112005
$ws.Cells.Item(224, 7).Value = "Puerro"
$ws.Cells.Item(224, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(224, 9).Value = "Primera"
$ws.Cells.Item(224, 10).Value = 45
$ws.Cells.Item(224, 11).Value = 14000
$ws.Cells.Item(224, 12).Value = 14000
$ws.Cells.Item(224, 13).Value = 14000
$ws.Cells.Item(224, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(224, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(224, 16).Value = 1167
$ws.Cells.Item(224, 17).Value = 12
$ws.Cells.Item(224, 18).Value = "Hortaliza"
